# Apply LDLC price-history update:
# Insert a new snapshot column "CX" (timestamp 2026-02-01 09:19:59), shifting
# the previous "nom" / "url_produit" columns one position to the right
# (CX->CY, CY->CZ). For rows 2-80 the new CX cell is populated with the same
# numeric price as column CW (last existing snapshot); for rows 81-206 (and
# beyond) the CW snapshot was empty so CX stays empty too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before CX; this shifts old CX -> CY and old CY -> CZ,
# carrying their values/formatting with them.
$ws.Columns("CX:CX").Insert()

# New header cell for the inserted snapshot column.
$ws.Range("CX1").Value = "2026-02-01 09:19:59"

# Populate the new snapshot column for every data row that still had a price
# in the previous snapshot (column CW), mirroring that same value.
$lastRow = 206
for ($r = 2; $r -le $lastRow; $r++) {
    $prev = $ws.Cells.Item($r, 101)   # column CW
    if ($prev.Value2 -ne "") {
        $cur = $ws.Cells.Item($r, 102)  # column CX (new)
        $cur.Value = $prev.Value2
    }
}
